$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
$v = $ws.Range("A1").Value
Write-Host ([string]$v)
Write-Host $ws.Range("A1").Value2
